$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 - team-members text box: drop "אסף עובדיה - 88888888" and the tab
# run that preceded "רווה פנחס - ", leaving a single run with the leading
# space kept (" רווה פנחס - ") directly before the "205357809" run.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$teamShape = $s1.Shapes.Item(2)
$teamRange = $teamShape.TextFrame.TextRange

$teamPara = $teamRange.Paragraphs(3)
$run4 = $teamPara.Runs(4)   # "88888888"
$run5 = $teamPara.Runs(5)   # tabs + "רווה פנחס - "

# Work out how many leading whitespace/tab characters run5 has, so only
# those are stripped (the trailing "רווה פנחס - " text must stay).
$run5Text = $run5.Text
$prefixLen = 0
while ($prefixLen -lt $run5Text.Length -and ($run5Text.Substring($prefixLen, 1) -eq "`t" -or $run5Text.Substring($prefixLen, 1) -eq " ")) {
    $prefixLen = $prefixLen + 1
}

$run4Start = $run4.Start
$run4Len = $run4.Text.Length
$run5Start = $run5.Start

# Delete right-to-left so earlier offsets remain valid.
if ($prefixLen -gt 0) {
    $teamRange.Characters($run5Start, $prefixLen).Delete()
}
$teamRange.Characters($run4Start, $run4Len).Delete()

# After the deletions, "רווה פנחס - " (formerly the tail of run5) became its
# own adjoining run right after run3 - merge them into run3 by deleting that
# leftover run and writing the combined text into run3.
$teamParaAfter = $teamRange.Paragraphs(3)
$mergedRun = $teamParaAfter.Runs(4)
$mergedText = $mergedRun.Text
$teamRange.Characters($mergedRun.Start, $mergedText.Length).Delete()

$teamParaFinal = $teamRange.Paragraphs(3)
$teamParaFinal.Runs(3).Text = " " + $mergedText

# ---------------------------------------------------------------------------
# Slide 8 - flow-chart: reposition two existing shapes, rename the last
# rectangle's label, then add a "rejection" branch (new rectangle + arrow).
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)

# "Straight Arrow Connector 14" -> move.
$conn14 = $s8.Shapes.Item(4)
$conn14.Left = 397.2135471217839
$conn14.Top = 211.8026809800516

# "Rectangle 5" -> move.
$rect5 = $s8.Shapes.Item(7)
$rect5.Left = 534.1334645669291
$rect5.Top = 136.4571653543307

# "Rectangle 7" -> rename label.
$rect7 = $s8.Shapes.Item(8)
$rect7.TextFrame.TextRange.Text = "אישור"

# New rectangle ("דחייה"), cloned from Rectangle 7 so it keeps the same
# theme-styled look (line/fill/effect/font refs) and text formatting.
$newRectSet = $rect7.Duplicate()
$newRect = $newRectSet.Item(1)
$newRect.Left = 402.19267335065706
$newRect.Top = 302.85054355170786
$newRect.TextFrame.TextRange.Text = "דחייה"

# New arrow connector, cloned from "Straight Arrow Connector 14" so the
# line weight / arrow head / theme colors match.
$newConnSet = $conn14.Duplicate()
$newConn = $newConnSet.Item(1)
$newConn.Left = 563.9772669826718
$newConn.Top = 252.27794703385962
$newConn.Width = 11.999999523162842
$newConn.Height = 62.49385683720506
